# edit.ps1
# Applies the commit: renames the "Requested quantity" headers on the two
# existing sheets and adds a new "PO Forecast" sheet (ds / PO_Forecast /
# yhat_lower / yhat_upper) with 49 rows of forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy the header formatting (bold, centered, bordered) from the Weekly
# Quantity sheet's header row so the new sheet's header matches the
# existing header style, then fill in the new header text.
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- 3. Fill in the forecast data rows (2-50) ---
$ws3.Range("A2").Value = 45032.99999999999
$ws3.Range("B2").Value = 387
$ws3.Range("C2").Value = -798.6293822339636
$ws3.Range("D2").Value = 1551.400806297615
$ws3.Range("A3").Value = 45046.99999999999
$ws3.Range("B3").Value = 418
$ws3.Range("C3").Value = -700.2021165632212
$ws3.Range("D3").Value = 1593.542619226176
$ws3.Range("A4").Value = 45053.99999999999
$ws3.Range("B4").Value = 434
$ws3.Range("C4").Value = -698.7770731242231
$ws3.Range("D4").Value = 1695.01683118264
$ws3.Range("A5").Value = 45060.99999999999
$ws3.Range("B5").Value = 449
$ws3.Range("C5").Value = -745.5169950668676
$ws3.Range("D5").Value = 1696.799713229344
$ws3.Range("A6").Value = 45067.99999999999
$ws3.Range("B6").Value = 465
$ws3.Range("C6").Value = -721.2198666148871
$ws3.Range("D6").Value = 1567.952316856707
$ws3.Range("A7").Value = 45074.99999999999
$ws3.Range("B7").Value = 481
$ws3.Range("C7").Value = -761.8876991351958
$ws3.Range("D7").Value = 1640.609052114338
$ws3.Range("A8").Value = 45088.99999999999
$ws3.Range("B8").Value = 512
$ws3.Range("C8").Value = -720.2692119819717
$ws3.Range("D8").Value = 1631.8493179123
$ws3.Range("A9").Value = 45095.99999999999
$ws3.Range("B9").Value = 528
$ws3.Range("C9").Value = -659.8596754552813
$ws3.Range("D9").Value = 1680.060909520509
$ws3.Range("A10").Value = 45102.99999999999
$ws3.Range("B10").Value = 543
$ws3.Range("C10").Value = -597.6258520057889
$ws3.Range("D10").Value = 1696.780504689848
$ws3.Range("A11").Value = 45109.99999999999
$ws3.Range("B11").Value = 559
$ws3.Range("C11").Value = -589.8019638023063
$ws3.Range("D11").Value = 1721.085149268527
$ws3.Range("A12").Value = 45123.99999999999
$ws3.Range("B12").Value = 590
$ws3.Range("C12").Value = -535.5977978940734
$ws3.Range("D12").Value = 1717.572959055443
$ws3.Range("A13").Value = 45137.99999999999
$ws3.Range("B13").Value = 622
$ws3.Range("C13").Value = -552.3352765625763
$ws3.Range("D13").Value = 1782.389930911799
$ws3.Range("A14").Value = 45144.99999999999
$ws3.Range("B14").Value = 637
$ws3.Range("C14").Value = -527.5371134771304
$ws3.Range("D14").Value = 1854.33259916792
$ws3.Range("A15").Value = 45151.99999999999
$ws3.Range("B15").Value = 653
$ws3.Range("C15").Value = -549.7928782102216
$ws3.Range("D15").Value = 1892.714336339296
$ws3.Range("A16").Value = 45172.99999999999
$ws3.Range("B16").Value = 700
$ws3.Range("C16").Value = -498.1770960548546
$ws3.Range("D16").Value = 1840.761390149699
$ws3.Range("A17").Value = 45179.99999999999
$ws3.Range("B17").Value = 716
$ws3.Range("C17").Value = -439.3413684617729
$ws3.Range("D17").Value = 1894.444409186503
$ws3.Range("A18").Value = 45200.99999999999
$ws3.Range("B18").Value = 763
$ws3.Range("C18").Value = -436.3558299955432
$ws3.Range("D18").Value = 1953.595628714323
$ws3.Range("A19").Value = 45207.99999999999
$ws3.Range("B19").Value = 778
$ws3.Range("C19").Value = -472.0487298613351
$ws3.Range("D19").Value = 1918.84688780884
$ws3.Range("A20").Value = 45221.99999999999
$ws3.Range("B20").Value = 810
$ws3.Range("C20").Value = -380.0397804217156
$ws3.Range("D20").Value = 2039.543423888914
$ws3.Range("A21").Value = 45228.99999999999
$ws3.Range("B21").Value = 825
$ws3.Range("C21").Value = -274.8470340821949
$ws3.Range("D21").Value = 2002.807131539005
$ws3.Range("A22").Value = 45263.99999999999
$ws3.Range("B22").Value = 904
$ws3.Range("C22").Value = -297.9893707824098
$ws3.Range("D22").Value = 2110.525780365101
$ws3.Range("A23").Value = 45270.99999999999
$ws3.Range("B23").Value = 919
$ws3.Range("C23").Value = -330.7529839421145
$ws3.Range("D23").Value = 2186.822045832304
$ws3.Range("A24").Value = 45277.99999999999
$ws3.Range("B24").Value = 935
$ws3.Range("C24").Value = -291.7640996144528
$ws3.Range("D24").Value = 2052.028887756309
$ws3.Range("A25").Value = 45298.99999999999
$ws3.Range("B25").Value = 982
$ws3.Range("C25").Value = -264.5399778220923
$ws3.Range("D25").Value = 2248.079922640864
$ws3.Range("A26").Value = 45312.99999999999
$ws3.Range("B26").Value = 1013
$ws3.Range("C26").Value = -167.2895076339827
$ws3.Range("D26").Value = 2153.783564261344
$ws3.Range("A27").Value = 45319.99999999999
$ws3.Range("B27").Value = 1029
$ws3.Range("C27").Value = -205.0344720615607
$ws3.Range("D27").Value = 2200.154585881031
$ws3.Range("A28").Value = 45326.99999999999
$ws3.Range("B28").Value = 1045
$ws3.Range("C28").Value = -87.20661232785304
$ws3.Range("D28").Value = 2222.767845681699
$ws3.Range("A29").Value = 45333.99999999999
$ws3.Range("B29").Value = 1060
$ws3.Range("C29").Value = -189.7301637737038
$ws3.Range("D29").Value = 2291.262906912455
$ws3.Range("A30").Value = 45340.99999999999
$ws3.Range("B30").Value = 1076
$ws3.Range("C30").Value = -140.5081595095018
$ws3.Range("D30").Value = 2189.041171659709
$ws3.Range("A31").Value = 45347.99999999999
$ws3.Range("B31").Value = 1092
$ws3.Range("C31").Value = -73.56573006288451
$ws3.Range("D31").Value = 2167.566010471062
$ws3.Range("A32").Value = 45417.99999999999
$ws3.Range("B32").Value = 1248
$ws3.Range("C32").Value = 58.35219753331198
$ws3.Range("D32").Value = 2479.529198646418
$ws3.Range("A33").Value = 45445.99999999999
$ws3.Range("B33").Value = 1310
$ws3.Range("C33").Value = 187.40071793926
$ws3.Range("D33").Value = 2556.941094270876
$ws3.Range("A34").Value = 45452.99999999999
$ws3.Range("B34").Value = 1326
$ws3.Range("C34").Value = 145.6540288378479
$ws3.Range("D34").Value = 2506.634151022301
$ws3.Range("A35").Value = 45459.99999999999
$ws3.Range("B35").Value = 1341
$ws3.Range("C35").Value = 192.7501080314793
$ws3.Range("D35").Value = 2488.129343979892
$ws3.Range("A36").Value = 45480.99999999999
$ws3.Range("B36").Value = 1388
$ws3.Range("C36").Value = 274.6299040233379
$ws3.Range("D36").Value = 2580.485654188221
$ws3.Range("A37").Value = 45487.99999999999
$ws3.Range("B37").Value = 1404
$ws3.Range("C37").Value = 190.5993138662529
$ws3.Range("D37").Value = 2621.639460482415
$ws3.Range("A38").Value = 45494.99999999999
$ws3.Range("B38").Value = 1419
$ws3.Range("C38").Value = 255.6076049487691
$ws3.Range("D38").Value = 2623.370875414723
$ws3.Range("A39").Value = 45501.99999999999
$ws3.Range("B39").Value = 1435
$ws3.Range("C39").Value = 278.2330286704495
$ws3.Range("D39").Value = 2594.952817284706
$ws3.Range("A40").Value = 45585.99999999999
$ws3.Range("B40").Value = 1622
$ws3.Range("C40").Value = 467.5923805782527
$ws3.Range("D40").Value = 2857.097620684582
$ws3.Range("A41").Value = 45599.99999999999
$ws3.Range("B41").Value = 1654
$ws3.Range("C41").Value = 464.2780325035035
$ws3.Range("D41").Value = 2865.551284928545
$ws3.Range("A42").Value = 45613.99999999999
$ws3.Range("B42").Value = 1685
$ws3.Range("C42").Value = 468.0321383331052
$ws3.Range("D42").Value = 2849.571089929602
$ws3.Range("A43").Value = 45620.99999999999
$ws3.Range("B43").Value = 1700
$ws3.Range("C43").Value = 489.9516229804909
$ws3.Range("D43").Value = 2951.494481923686
$ws3.Range("A44").Value = 45627.99999999999
$ws3.Range("B44").Value = 1716
$ws3.Range("C44").Value = 555.8234793178901
$ws3.Range("D44").Value = 2854.569590156489
$ws3.Range("A45").Value = 45634.99999999999
$ws3.Range("B45").Value = 1732
$ws3.Range("C45").Value = 472.6633471411058
$ws3.Range("D45").Value = 2944.281590167745
$ws3.Range("A46").Value = 45641.99999999999
$ws3.Range("B46").Value = 1747
$ws3.Range("C46").Value = 615.501036983047
$ws3.Range("D46").Value = 2902.066741142262
$ws3.Range("A47").Value = 45648.99999999999
$ws3.Range("B47").Value = 1763
$ws3.Range("C47").Value = 561.5449850820019
$ws3.Range("D47").Value = 2911.122980101627
$ws3.Range("A48").Value = 45655.99999999999
$ws3.Range("B48").Value = 1779
$ws3.Range("C48").Value = 567.9100487896283
$ws3.Range("D48").Value = 2889.664880704264
$ws3.Range("A49").Value = 45662.99999999999
$ws3.Range("B49").Value = 1794
$ws3.Range("C49").Value = 576.4865395978796
$ws3.Range("D49").Value = 2993.936559037046
$ws3.Range("A50").Value = 45669.99999999999
$ws3.Range("B50").Value = 1810
$ws3.Range("C50").Value = 591.3628818249186
$ws3.Range("D50").Value = 2964.561149508793

# Match the date-time number format already used for column A on the other
# sheets (style reused automatically since the format code is identical).
$ws3.Range("A2:A50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
